# Refresh the benchmark numbers on Sheet1 (three Sequential-vs-parallel
# tables) to match the latest timing run, and tidy up the view that was
# left scrolled down/selected on a stray cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Table 1 (rows 4-7, header row 2): Sequential (B) vs Loop-level parallelism (C) ---
$ws.Range("B4").Value = 366
$ws.Range("C4").Value = 204
$ws.Range("B5").Value = 238
$ws.Range("C5").Value = 117
$ws.Range("B6").Value = 239
$ws.Range("C6").Value = 124
$ws.Range("B7").Value = 241
$ws.Range("C7").Value = 157

# Column B in this table now matches column A/C's plain centred style
# instead of the 2-decimal-place style it had before. Re-use the existing
# style via a formats-only paste so no new style entries get created.
$ws.Range("A4").Copy()
$ws.Range("B4:B7").PasteSpecial(-4122)

# --- Table 2 (rows 14-17, header row 12): Sequential (B) vs functional decomposition (C) ---
$ws.Range("B14").Value = 366
$ws.Range("C14").Value = 248
$ws.Range("B15").Value = 238
$ws.Range("C15").Value = 215
$ws.Range("B16").Value = 239
$ws.Range("C16").Value = 232
$ws.Range("B17").Value = 241
$ws.Range("C17").Value = 233

$ws.Range("A14").Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)

# --- Table 3 (rows 24-27, header row 22): Sequential (B), functional decomposition (C), Loop-level (D) ---
$ws.Range("B24").Value = 366
$ws.Range("C24").Value = 248
$ws.Range("D24").Value = 204
$ws.Range("B25").Value = 238
$ws.Range("C25").Value = 215
$ws.Range("B26").Value = 239
$ws.Range("C26").Value = 232
$ws.Range("D26").Value = 124
$ws.Range("B27").Value = 241
$ws.Range("C27").Value = 233
$ws.Range("D27").Value = 157

$excel.CutCopyMode = 0

# --- Restore the view: scroll back to the top and move the selection ---
$ws.Activate()
$ws.Range("E32").Select() | Out-Null
